$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before row 998, shifting the existing data
# (rows 998-1078) down to rows 1000-1080.
$ws.Rows("998:999").Insert()

# New row 998: Primera quality entry for the new week
$ws.Cells.Item(998, 1).Value2  = 9
$ws.Cells.Item(998, 2).Value2  = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(998, 3).Value2  = "Metropolitana"
$ws.Cells.Item(998, 4).Value2  = 45223
$ws.Cells.Item(998, 5).Value2  = 13
$ws.Cells.Item(998, 6).Value2  = 100112008
$ws.Cells.Item(998, 7).Value2  = "Coliflor"
$ws.Cells.Item(998, 8).Value2  = "Sin especificar"
$ws.Cells.Item(998, 9).Value2  = "Primera"
$ws.Cells.Item(998, 10).Value2 = 1600
$ws.Cells.Item(998, 11).Value2 = 800
$ws.Cells.Item(998, 12).Value2 = 900
$ws.Cells.Item(998, 13).Value2 = 850
$ws.Cells.Item(998, 14).Value2 = "`$/unidad"
$ws.Cells.Item(998, 15).Value2 = "Región Metropolitana"
$ws.Cells.Item(998, 16).Value2 = 850
$ws.Cells.Item(998, 17).Value2 = 1
$ws.Cells.Item(998, 18).Value2 = "Hortaliza"

# New row 999: Segunda quality entry for the new week
$ws.Cells.Item(999, 1).Value2  = 9
$ws.Cells.Item(999, 2).Value2  = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(999, 3).Value2  = "Metropolitana"
$ws.Cells.Item(999, 4).Value2  = 45223
$ws.Cells.Item(999, 5).Value2  = 13
$ws.Cells.Item(999, 6).Value2  = 100112008
$ws.Cells.Item(999, 7).Value2  = "Coliflor"
$ws.Cells.Item(999, 8).Value2  = "Sin especificar"
$ws.Cells.Item(999, 9).Value2  = "Segunda"
$ws.Cells.Item(999, 10).Value2 = 970
$ws.Cells.Item(999, 11).Value2 = 700
$ws.Cells.Item(999, 12).Value2 = 700
$ws.Cells.Item(999, 13).Value2 = 700
$ws.Cells.Item(999, 14).Value2 = "`$/unidad"
$ws.Cells.Item(999, 15).Value2 = "Región Metropolitana"
$ws.Cells.Item(999, 16).Value2 = 700
$ws.Cells.Item(999, 17).Value2 = 1
$ws.Cells.Item(999, 18).Value2 = "Hortaliza"
